$wb = $excel.ActiveWorkbook

# --- Sheet: LP1912 ---
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2, 1).Value = 'Última actualización: 20:32:11'
$ws.Cells.Item(3, 1).Value = 'Total filas: 362'
$ws.Cells.Item(38, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(39, 3).Value = '15_ABASTO'
$ws.Cells.Item(49, 1).Value = '07:13:03'
$ws.Cells.Item(49, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(49, 4).Value = 99
$ws.Cells.Item(50, 1).Value = '08:52:40'
$ws.Cells.Item(50, 3).Value = '215B_EL PATO'
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(94, 1).Value = '10:56:15'
$ws.Cells.Item(94, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(94, 4).Value = 1
$ws.Cells.Item(95, 1).Value = '09:22:34'
$ws.Cells.Item(95, 3).Value = '10_OLMOS'
$ws.Cells.Item(95, 4).Value = 95
$ws.Cells.Item(139, 1).Value = '11:53:44'
$ws.Cells.Item(139, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(139, 4).Value = 43
$ws.Cells.Item(140, 1).Value = '10:49:38'
$ws.Cells.Item(140, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(140, 4).Value = 107
$ws.Cells.Item(276, 1).Value = '16:37:37'
$ws.Cells.Item(276, 3).Value = '14_ABASTO'
$ws.Cells.Item(276, 4).Value = 112
$ws.Cells.Item(277, 1).Value = '17:13:30'
$ws.Cells.Item(277, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(277, 4).Value = 76
$ws.Cells.Item(334, 1).Value = '20:32:11'
$ws.Cells.Item(334, 2).Value = '20:38'
$ws.Cells.Item(334, 3).Value = '10_OLMOS'
$ws.Cells.Item(334, 4).Value = 6
$ws.Cells.Item(335, 1).Value = '19:11:44'
$ws.Cells.Item(335, 2).Value = '20:41'
$ws.Cells.Item(335, 4).Value = 90
$ws.Cells.Item(336, 1).Value = '18:52:29'
$ws.Cells.Item(336, 2).Value = '20:42'
$ws.Cells.Item(336, 4).Value = 110
$ws.Cells.Item(337, 1).Value = '19:35:34'
$ws.Cells.Item(337, 2).Value = '20:43'
$ws.Cells.Item(337, 4).Value = 68
$ws.Cells.Item(338, 1).Value = '19:47:50'
$ws.Cells.Item(338, 2).Value = '20:45'
$ws.Cells.Item(338, 3).Value = '17_ROMERO'
$ws.Cells.Item(338, 4).Value = 58
$ws.Cells.Item(339, 1).Value = '20:32:11'
$ws.Cells.Item(339, 2).Value = '20:46'
$ws.Cells.Item(339, 3).Value = '17_ROMERO'
$ws.Cells.Item(339, 4).Value = 14
$ws.Cells.Item(340, 1).Value = '18:52:29'
$ws.Cells.Item(340, 2).Value = '20:47'
$ws.Cells.Item(340, 3).Value = '215B_EL PATO'
$ws.Cells.Item(340, 4).Value = 115
$ws.Cells.Item(341, 1).Value = '20:32:11'
$ws.Cells.Item(341, 2).Value = '20:48'
$ws.Cells.Item(341, 3).Value = '215B_EL PATO'
$ws.Cells.Item(341, 4).Value = 16
$ws.Cells.Item(342, 1).Value = '20:11:58'
$ws.Cells.Item(342, 2).Value = '20:54'
$ws.Cells.Item(342, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(342, 4).Value = 43
$ws.Cells.Item(343, 1).Value = '19:35:34'
$ws.Cells.Item(343, 2).Value = '20:55'
$ws.Cells.Item(343, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(343, 4).Value = 80
$ws.Cells.Item(344, 1).Value = '19:54:57'
$ws.Cells.Item(344, 2).Value = '20:55'
$ws.Cells.Item(344, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(344, 4).Value = 61
$ws.Cells.Item(345, 1).Value = '19:11:44'
$ws.Cells.Item(345, 2).Value = '20:56'
$ws.Cells.Item(345, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(345, 4).Value = 105
$ws.Cells.Item(346, 1).Value = '20:32:11'
$ws.Cells.Item(346, 2).Value = '20:57'
$ws.Cells.Item(346, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(346, 4).Value = 25
$ws.Cells.Item(347, 2).Value = '21:06'
$ws.Cells.Item(347, 3).Value = '14_ABASTO'
$ws.Cells.Item(347, 4).Value = 55
$ws.Cells.Item(348, 1).Value = '19:11:44'
$ws.Cells.Item(348, 2).Value = '21:06'
$ws.Cells.Item(348, 3).Value = '10_OLMOS'
$ws.Cells.Item(348, 4).Value = 115
$ws.Cells.Item(349, 1).Value = '20:32:11'
$ws.Cells.Item(349, 2).Value = '21:07'
$ws.Cells.Item(349, 3).Value = '10_OLMOS'
$ws.Cells.Item(349, 4).Value = 35
$ws.Cells.Item(350, 2).Value = '21:09'
$ws.Cells.Item(350, 3).Value = '15_ABASTO'
$ws.Cells.Item(350, 4).Value = 82
$ws.Cells.Item(351, 1).Value = '20:32:11'
$ws.Cells.Item(351, 2).Value = '21:09'
$ws.Cells.Item(351, 3).Value = '14_ABASTO'
$ws.Cells.Item(351, 4).Value = 37
$ws.Cells.Item(352, 1).Value = '19:35:34'
$ws.Cells.Item(352, 2).Value = '21:10'
$ws.Cells.Item(352, 3).Value = '15_ABASTO'
$ws.Cells.Item(352, 4).Value = 95
$ws.Cells.Item(353, 1).Value = '20:11:58'
$ws.Cells.Item(353, 2).Value = '21:27'
$ws.Cells.Item(353, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(353, 4).Value = 76
$ws.Cells.Item(354, 1).Value = '19:35:34'
$ws.Cells.Item(354, 2).Value = '21:28'
$ws.Cells.Item(354, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(354, 4).Value = 113
$ws.Cells.Item(355, 1).Value = '19:54:57'
$ws.Cells.Item(355, 2).Value = '21:33'
$ws.Cells.Item(355, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(355, 4).Value = 99
$ws.Cells.Item(356, 1).Value = '19:47:50'
$ws.Cells.Item(356, 2).Value = '21:33'
$ws.Cells.Item(356, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(356, 4).Value = 106
$ws.Cells.Item(357, 1).Value = '19:35:34'
$ws.Cells.Item(357, 2).Value = '21:34'
$ws.Cells.Item(357, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(357, 4).Value = 119
$ws.Cells.Item(357, 5).Value = 'LP1912'
$ws.Cells.Item(358, 1).Value = '20:11:58'
$ws.Cells.Item(358, 2).Value = '21:37'
$ws.Cells.Item(358, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(358, 4).Value = 86
$ws.Cells.Item(358, 5).Value = 'LP1912'
$ws.Cells.Item(359, 1).Value = '20:32:11'
$ws.Cells.Item(359, 2).Value = '21:39'
$ws.Cells.Item(359, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(359, 4).Value = 67
$ws.Cells.Item(359, 5).Value = 'LP1912'
$ws.Cells.Item(360, 1).Value = '19:54:57'
$ws.Cells.Item(360, 2).Value = '21:44'
$ws.Cells.Item(360, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(360, 4).Value = 110
$ws.Cells.Item(360, 5).Value = 'LP1912'
$ws.Cells.Item(361, 1).Value = '19:47:50'
$ws.Cells.Item(361, 2).Value = '21:45'
$ws.Cells.Item(361, 3).Value = '14X44_ABASTO'
$ws.Cells.Item(361, 4).Value = 118
$ws.Cells.Item(361, 5).Value = 'LP1912'
$ws.Cells.Item(362, 1).Value = '20:32:11'
$ws.Cells.Item(362, 2).Value = '21:46'
$ws.Cells.Item(362, 3).Value = '14X44_ABASTO'
$ws.Cells.Item(362, 4).Value = 74
$ws.Cells.Item(362, 5).Value = 'LP1912'
$ws.Cells.Item(363, 1).Value = '20:11:58'
$ws.Cells.Item(363, 2).Value = '21:48'
$ws.Cells.Item(363, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(363, 4).Value = 97
$ws.Cells.Item(363, 5).Value = 'LP1912'
$ws.Cells.Item(364, 1).Value = '20:32:11'
$ws.Cells.Item(364, 2).Value = '21:51'
$ws.Cells.Item(364, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(364, 4).Value = 79
$ws.Cells.Item(364, 5).Value = 'LP1912'
$ws.Cells.Item(365, 1).Value = '20:11:58'
$ws.Cells.Item(365, 2).Value = '22:03'
$ws.Cells.Item(365, 3).Value = '15_ABASTO'
$ws.Cells.Item(365, 4).Value = 112
$ws.Cells.Item(365, 5).Value = 'LP1912'
$ws.Cells.Item(366, 1).Value = '20:32:11'
$ws.Cells.Item(366, 2).Value = '22:04'
$ws.Cells.Item(366, 3).Value = '15_ABASTO'
$ws.Cells.Item(366, 4).Value = 92
$ws.Cells.Item(366, 5).Value = 'LP1912'
$ws.Cells.Item(367, 1).Value = '20:32:11'
$ws.Cells.Item(367, 2).Value = '22:11'
$ws.Cells.Item(367, 3).Value = '14_ABASTO'
$ws.Cells.Item(367, 4).Value = 99
$ws.Cells.Item(367, 5).Value = 'LP1912'

# --- Sheet: LP1912-215 ---
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2, 1).Value = 'Última actualización: 20:32:11'
$ws.Cells.Item(3, 1).Value = 'Total filas: 54'
$ws.Cells.Item(59, 1).Value = '20:32:11'
$ws.Cells.Item(59, 2).Value = '20:48'
$ws.Cells.Item(59, 3).Value = '215B_EL PATO'
$ws.Cells.Item(59, 4).Value = 16
$ws.Cells.Item(59, 5).Value = 'LP1912'

# --- Sheet: 6203-6173 ---
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2, 1).Value = 'Última actualización: 20:32:11'
$ws.Cells.Item(3, 1).Value = 'Total filas: 51'
$ws.Cells.Item(52, 1).Value = '20:32:11'
$ws.Cells.Item(52, 2).Value = '20:52'
$ws.Cells.Item(52, 4).Value = 20
$ws.Cells.Item(53, 1).Value = '19:47:50'
$ws.Cells.Item(53, 2).Value = '21:27'
$ws.Cells.Item(53, 4).Value = 100
$ws.Cells.Item(54, 1).Value = '19:54:57'
$ws.Cells.Item(54, 2).Value = '21:29'
$ws.Cells.Item(54, 4).Value = 95
$ws.Cells.Item(55, 1).Value = '19:35:34'
$ws.Cells.Item(55, 2).Value = '21:30'
$ws.Cells.Item(55, 3).Value = '215C_LA PLATA'
$ws.Cells.Item(55, 4).Value = 115
$ws.Cells.Item(55, 5).Value = 'L6203'
$ws.Cells.Item(56, 1).Value = '20:32:11'
$ws.Cells.Item(56, 2).Value = '22:20'
$ws.Cells.Item(56, 3).Value = '215B_LP-P MOR-40 Y 115'
$ws.Cells.Item(56, 4).Value = 108
$ws.Cells.Item(56, 5).Value = 'L6173'
